$d = $word.ActiveDocument

# NOTE on ordering: in this runtime, touching the Tables collection
# before walking $d.Paragraphs causes the Paragraphs collection's
# cached Range offsets to go stale. So we resolve and delete the
# blank paragraphs (step 1, pure Paragraphs/Range work) BEFORE we
# touch $d.Tables at all (step 2).

# ------------------------------------------------------------------
# 1) Right after the UKL/UPL verification table there are four blank
#    paragraphs sharing identical formatting (justified, right indent
#    -606 twips / -30.3pt) followed by one differently formatted
#    paragraph. Remove the first two of those four duplicates, using
#    that formatting fingerprint (rather than a hard-coded paragraph
#    index) so the match is unambiguous and unique in the document.
# ------------------------------------------------------------------

$removed = 0
$guard = 0
while ($removed -lt 2 -and $guard -lt 10) {
    $guard = $guard + 1
    $cnt = $d.Paragraphs.Count
    $hit = -1
    for ($i = 1; $i -le $cnt; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Trim() -eq "" -and $p.Alignment -eq 3) {
            $ri = [Math]::Round($p.Range.ParagraphFormat.RightIndent, 1)
            if ($ri -eq -30.3) {
                $hit = $i
                break
            }
        }
    }
    if ($hit -eq -1) {
        break
    }
    $d.Paragraphs.Item($hit).Range.Delete()
    $removed = $removed + 1
}

# ------------------------------------------------------------------
# 2) Delete the table row "11 - Penambahan Peta Titk Pengelolaan dan
#    Titik Pemantauan" (the ${peta_titik_*} placeholders row) from
#    the verification table.
# ------------------------------------------------------------------

for ($ti = 1; $ti -le $d.Tables.Count; $ti++) {
    $table = $d.Tables.Item($ti)
    for ($r = $table.Rows.Count; $r -ge 1; $r--) {
        $row = $table.Rows.Item($r)
        if ($row.Cells.Item(2).Range.Text -like "*Penambahan Peta Titk Pengelolaan dan Titik Pemantauan*") {
            $row.Delete()
            break
        }
    }
}
